# Generate Report for Handoff
#
# Re-running the handoff-generation step produced a fresh Xliff handoff
# batch for the "6a186769-89f5-4b00-9ebf-497fdf37052a.md" source file and
# everything that depends on it (rows 4-7 of the per-locale sheets).
# Its Priority flips from "low" to "ht" (high priority / hot?) and the
# "Latest Handoff Datetime" for that batch is refreshed to the new
# generation timestamp, for both locales (zh-cn and de-de).

$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
$zh.Range("E4:E7").Value = "ht"
$zh.Range("H4:H7").Value = "2016-08-12 06:37:13"

# de-de: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
$de.Range("E4:E7").Value = "ht"
$de.Range("H4:H7").Value = "2016-08-12 06:37:20"

# Overview: rows 4-7 -> Latest HO Xliff Generate Date (G) mirrors the de-de
# handoff batch timestamp, so it is refreshed too.
$ov.Range("G4:G7").Value = "2016-08-12 06:37:20"
